# تعديل تلقائي في شيت Card2 by admin at 2025-11-02 09:29:54
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# Header O1: drop the trailing space ("Serviced by " -> "Serviced by")
$ws.Range("O1").Value = "Serviced by"

# Fill the previously-empty "Serviced by" column (O) with "nan" placeholders,
# matching the rest of the row's empty-value convention.
$ws.Range("O2").Value = "nan"
$ws.Range("O3").Value = "nan"
$ws.Range("O4").Value = "nan"
$ws.Range("O5").Value = "nan"
$ws.Range("O6").Value = "nan"
$ws.Range("O7").Value = "nan"
$ws.Range("O8").Value = "nan"

# Row 9 gets real service data instead of placeholders.
$ws.Range("M9").Value = "تم سن فلاتس وعياره "
$ws.Range("O9").Value = "م.محمد عبدالله "

$ws.Range("O10").Value = "nan"
$ws.Range("O11").Value = "nan"
$ws.Range("O12").Value = "nan"
$ws.Range("O13").Value = "nan"
